$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.354.48'
$ws.Range("E2").Value = '  +0.59%  '
$ws.Range("D3").Value = '2.628.86'
$ws.Range("E3").Value = '  +0.77%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''601.90'
$ws.Range("E5").Value = '  +1.73%  '
$ws.Range("D6").Value = '''153.43'
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '''0.557'
$ws.Range("E8").Value = '  +2.89%  '
$ws.Range("D9").Value = '2.626.00'
$ws.Range("E9").Value = '  +0.69%  '
$ws.Range("E10").Value = '  +6.13%  '
$ws.Range("E11").Value = '  +0.68%  '
$ws.Range("D12").Value = '''5.20'
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("E13").Value = '  -0.57%  '
$ws.Range("D14").Value = '''27.92'
$ws.Range("E14").Value = '  +1.36%  '
$ws.Range("D15").Value = '3.113.42'
$ws.Range("E15").Value = '  +1.02%  '
$ws.Range("D16").Value = '''0.0000183'
$ws.Range("E16").Value = '  +1.15%  '
$ws.Range("D17").Value = '67.290.73'
$ws.Range("E17").Value = '  +0.39%  '
$ws.Range("D18").Value = '2.630.64'
$ws.Range("E18").Value = '  +0.76%  '
$ws.Range("D19").Value = '''11.24'
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").Value = '''364.33'
$ws.Range("E20").Value = '  +1.44%  '
$ws.Range("D21").Value = '''7.60'
$ws.Range("E21").Value = '  -3.64%  '
$ws.Range("D22").Value = '''4.30'
$ws.Range("E22").Value = '  -0.33%  '
$ws.Range("D23").Value = '''2.13'
$ws.Range("E23").Value = '  +6.03%  '
$ws.Range("D24").Value = '''1.00'
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").Value = '''10.11'
$ws.Range("E25").Value = '  -0.38%  '
$ws.Range("D26").Value = '''67.36'
$ws.Range("E26").Value = '  -6.16%  '
$ws.Range("D27").Value = '2.769.81'
$ws.Range("E28").Value = '  +0.88%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("B30").Value = 'Bittensor'
$ws.Range("C30").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D30").Value = '''579.45'
$ws.Range("E30").Value = '  -6.41%  '
$ws.Range("D31").Value = '''1.41'
$ws.Range("E31").Value = '  -2.37%  '
$ws.Range("D32").Value = '''7.89'
$ws.Range("D33").Value = '''1.85'
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D34").Value = '''1.00'
$ws.Range("E34").Value = '  +0.13%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").Value = '''0.128'
$ws.Range("E35").Value = '  -3.59%  '
$ws.Range("D36").Value = '''1.53'
$ws.Range("D37").Value = '''4.94'
$ws.Range("E37").Value = '  -0.49%  '
$ws.Range("D38").Value = '''158.31'
$ws.Range("E38").Value = '  +3.06%  '
$ws.Range("D39").Value = '''19.40'
$ws.Range("E39").Value = '  +0.38%  '
$ws.Range("D40").Value = '''0.370'
$ws.Range("E40").Value = '  +0.52%  '
$ws.Range("D41").Value = '''5.29'
$ws.Range("E41").Value = '  -3.20%  '
$ws.Range("D42").Value = '''1.82'
$ws.Range("E42").Value = '  +0.33%  '
$ws.Range("E43").Value = '  +1.80%  '
$ws.Range("D44").Value = '''41.21'
$ws.Range("E44").Value = '  -0.32%  '
$ws.Range("E46").Value = '  -0.70%  '
$ws.Range("D47").Value = '''156.12'
$ws.Range("D48").Value = '0.0₆0288'
$ws.Range("E48").Value = '  -2.71%  '
$ws.Range("D49").Value = '''3.73'
$ws.Range("E49").Value = '  -0.76%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '''20.91'
$ws.Range("E50").Value = '  +0.33%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '''0.625'
$ws.Range("E51").Value = '  +0.53%  '
